$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BH: "Agosto.2021" period -------------------------------
# Copy the header cell's formatting (bold, centered, bordered) from BG1
# onto BH1, then write the new header label.
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)
$ws.Range("BH1").Value = "Agosto.2021"

# For every existing cohort row (2-73) the newest period simply repeats
# the prior period's value (carried forward), same as column BG.
$src = $ws.Range("BG2:BG73")
$dst = $ws.Range("BH2:BH73")
$dst.Value = $src.Value()

# Row 74 (cohort started Mayo.2021) gets its own new reading for
# Agosto.2021, distinct from its BG value.
$ws.Range("BH74").Value = 34514

# --- New row 75: cohort starting 01-04-2021 -----------------------------
# Build the label as a text string (avoiding Excel's automatic date
# coercion) via a helper formula cell, then paste the computed value back
# so no new number-format/style gets attached to A75.
$ws.Range("ZZ500").Formula = "=""01-04-2021"""
$ws.Range("ZZ500").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$ws.Range("ZZ500").Clear()

$ws.Range("BH75").Value = 35937
